# Update countries & provincias Spain
# - Reorders two country-name entries ("Libia" up next to "Yemen",
#   "Groenlandia" up next to "San Cristobal y Nieves") which, combined with
#   the underlying data refresh/re-sort, shifts the statistics that land on
#   a handful of rows.
# - Refreshes the "Datos actualizados ..." timestamp string.
# - Refreshes numeric COVID case figures for the rows that changed between
#   the two data pulls.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country names that moved position in the underlying data pull -------
# Yemen (row 128) and Jordania (row 133) keep their figures; the four rows
# in between now carry the next country's data because "Libia" sorted in
# right after "Yemen".
$ws.Range("A129").Value = "Libia"
$ws.Range("A130").Value = "Tunez"
$ws.Range("A131").Value = "Benin"
$ws.Range("A132").Value = "Ruanda"

# "Groenlandia" sorted in right after "San Cristobal y Nieves", swapping
# places with "Islas Malvinas" (their figures are identical, so only the
# labels change).
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Timestamp header -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Julio de 2020 a las 02:52"

# --- Refreshed statistics ----------------------------------------------
# Estados Unidos (row 4)
$ws.Range("B4").Value = 3158864
$ws.Range("C4").Value = 61780
$ws.Range("D4").Value = 1392679
$ws.Range("E4").Value = 1631325
$ws.Range("G4").Value = 888
$ws.Range("H4").Value = 134860

# Canada (row 23)
$ws.Range("B23").Value = 106434
$ws.Range("C23").Value = 267
$ws.Range("D23").Value = 70247
$ws.Range("E23").Value = 27450
$ws.Range("G23").Value = 26
$ws.Range("H23").Value = 8737

# Venezuela (row 79)
$ws.Range("E79").Value = 5833
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 75

# Libia (row 129, see rename above)
$ws.Range("B129").Value = 1268
$ws.Range("C129").Value = 86
$ws.Range("D129").Value = 306
$ws.Range("E129").Value = 926
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 36

# Tunez (row 130, see rename above)
$ws.Range("B130").Value = 1221
$ws.Range("C130").Value = 16
$ws.Range("D130").Value = 1050
$ws.Range("E130").Value = 121
$ws.Range("H130").Value = 50

# Benin (row 131, see rename above)
$ws.Range("B131").Value = 1199
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 333
$ws.Range("E131").Value = 845
$ws.Range("H131").Value = 21

# Ruanda (row 132, see rename above)
$ws.Range("B132").Value = 1194
$ws.Range("C132").Value = 22
$ws.Range("D132").Value = 610
$ws.Range("E132").Value = 581
$ws.Range("H132").Value = 3

# Niger (row 136)
$ws.Range("B136").Value = 1097
$ws.Range("C136").Value = 3
$ws.Range("D136").Value = 976
$ws.Range("E136").Value = 53

# Principado de Andorra (row 147)
$ws.Range("D147").Value = 802
$ws.Range("E147").Value = 1

# Togo (row 152)
$ws.Range("B152").Value = 695
$ws.Range("C152").Value = 6
$ws.Range("D152").Value = 475
$ws.Range("E152").Value = 205

# Islas Caimanes (row 171)
$ws.Range("D171").Value = 197
$ws.Range("E171").Value = 3

# Bermudas (row 176)
$ws.Range("B176").Value = 149
$ws.Range("C176").Value = 1
$ws.Range("E176").Value = 3
